$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("D1 Project")
$ws2 = $wb.Worksheets.Item("D2 Project")

# Clear the erroneous Input value for the first data row on both sheets
$ws1.Range("B5").ClearContents()
$ws2.Range("B5").ClearContents()

# Leave selection on B5 for both sheets, with the second sheet active/selected
$ws1.Range("B5").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B5").Select() | Out-Null
